$wb = $excel.ActiveWorkbook

# Update the two shared-string-backed values on the O_Metric sheet
$ws2 = $wb.Worksheets.Item("O_Metric")
$ws2.Range("B3").Value = "INTEREST"
$ws2.Range("B4").Value = "PRINCIPAL"

# Update selection on O_Metric sheet to B4 and make it the active sheet/tab
$ws2.Activate()
$ws2.Range("B4").Select()

$wb.Save()
